# ---------------------------------------------------------------------------
# Add the "Solution" worksheet to the water-consumption descriptive
# statistics exercise: quartiles, IQR, fences, outlier counts and skewness,
# all driven off a new "water" defined name.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item(1)

# Defined name covering the data column (A2:A351 on "water consumption").
$wb.Names.Add("water", "='water consumption'!`$A`$2:`$A`$351")

# New worksheet placed after the data sheet.
$ws = $wb.Worksheets.Add()
$ws.Name = "Solution"
$ws.Move($null, $wsData)

# ---- Header row --------------------------------------------------------------
$ws.Range("A1").Value = "Statistic"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Units"
$ws.Range("D1").Value = "Interpretation"

$hdr = $ws.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.Font.Name = "Arial"
$hdr.Font.Size = 10

# ---- Quartiles section ---------------------------------------------------------
$ws.Range("A2").Value = "Quartiles"
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Size = 10

$ws.Range("A3").Value = 1
$ws.Range("B3").Formula = "=QUARTILE(water,A3)"
$ws.Range("D3").Value = "25% of the househods consumed a quantity of water lower than or equal to this quantity. "
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D3").Font.Size = 10

$ws.Range("A4").Value = 2
$ws.Range("B4").Formula = "=QUARTILE(water,A4)"
$ws.Range("D4").Value = "50% of the househods consumed a quantity of water lower than or equal to this quantity. "
$ws.Range("D4").Font.Name = "Arial"
$ws.Range("D4").Font.Size = 10

$ws.Range("A5").Value = 3
$ws.Range("B5").Formula = "=QUARTILE(water,A5)"
$ws.Range("D5").Value = "75% of the househods consumed a quantity of water lower than or equal to this quantity. "
$ws.Range("D5").Font.Name = "Arial"
$ws.Range("D5").Font.Size = 10

# Build the rich "m3" (superscript 3) unit label once, then fan it out by copy
# so the three cells share a single rich-text entry, matching how Excel
# itself dedupes shared strings when a formatted cell is copy/pasted.
$ws.Range("C3").Value = "'m3"
$ws.Range("C3").Characters(2, 1).Font.Superscript = $true
$ws.Range("C3").Font.Name = "Arial"
$ws.Range("C3").Font.Size = 10
$ws.Range("C3").Copy($ws.Range("C4")) | Out-Null
$ws.Range("C3").Copy($ws.Range("C5")) | Out-Null

# ---- IQR and fences ----------------------------------------------------------
$ws.Range("A6").Value = "IQR"
$ws.Range("A6").Font.Name = "Arial"
$ws.Range("A6").Font.Size = 10
$ws.Range("B6").Formula = "=B5-B3"
$ws.Range("D6").Value = "Measures the spread of the 50% central data."
$ws.Range("D6").Font.Name = "Arial"
$ws.Range("D6").Font.Size = 10

$ws.Range("A7").Value = "lower fence"
$ws.Range("A7").Font.Name = "Arial"
$ws.Range("A7").Font.Size = 10
$ws.Range("B7").Formula = "=B3-1.5*B6"
$ws.Range("D7").Value = "Values lower than this value are outliers."
$ws.Range("D7").Font.Name = "Arial"
$ws.Range("D7").Font.Size = 10

$ws.Range("A8").Value = "upper fence"
$ws.Range("A8").Font.Name = "Arial"
$ws.Range("A8").Font.Size = 10
$ws.Range("B8").Formula = "=B5+1.5*B6"
$ws.Range("D8").Value = "Values greater than this value are outliers."
$ws.Range("D8").Font.Name = "Arial"
$ws.Range("D8").Font.Size = 10

# Second "m3" unit label (IQR / fences group) - its own shared, copy-fanned cell.
$ws.Range("C6").Value = "'m3"
$ws.Range("C6").Characters(2, 1).Font.Superscript = $true
$ws.Range("C6").Font.Name = "Arial"
$ws.Range("C6").Font.Size = 10
$ws.Range("C6").Copy($ws.Range("C7")) | Out-Null
$ws.Range("C6").Copy($ws.Range("C8")) | Out-Null

# ---- Outliers -----------------------------------------------------------------
$ws.Range("A9").Value = "Lower outliers"
$ws.Range("A9").Font.Name = "Arial"
$ws.Range("A9").Font.Size = 10
$ws.Range("B9").Formula = '=COUNTIF(water,"<"&B7)'

$ws.Range("A10").Value = "Upper outliers"
$ws.Range("A10").Font.Name = "Arial"
$ws.Range("A10").Font.Size = 10
$ws.Range("B10").Formula = '=COUNTIF(water,">"&B8)'

$ws.Range("A11").Value = "Total outliers"
$ws.Range("A11").Font.Name = "Arial"
$ws.Range("A11").Font.Size = 10
$ws.Range("B11").Formula = "=SUM(B9:B10)"

# ---- Skewness -------------------------------------------------------------------
$ws.Range("A12").Value = "Coef.Skewness"
$ws.Range("A12").Font.Name = "Arial"
$ws.Range("A12").Font.Size = 10
$ws.Range("B12").Formula = "=SKEW(water)"
$ws.Range("D12").Value = "As this value is positive the distribution is right-skewed. That means that there are a lot of households with low water consumption and few households with huge water consumption."
$ws.Range("D12").Font.Name = "Arial"
$ws.Range("D12").Font.Size = 10

# ---- Column widths (approximate bestFit sizing from the source workbook) -------
$ws.Columns.Item(1).ColumnWidth = 14
$ws.Columns.Item(2).ColumnWidth = 11
$ws.Columns.Item(3).ColumnWidth = 5.43

# ---- Selection / active sheet bookkeeping --------------------------------------
$ws.Range("D13").Select() | Out-Null
$wsData.Activate() | Out-Null

Write-Host "Solution sheet added"
